$wb = $excel.ActiveWorkbook

$wsScenarios = $wb.Worksheets.Item("Test Scenarios")
$wsCases     = $wb.Worksheets.Item("Test Cases")
$wsSheet3    = $wb.Worksheets.Item("Sheet3")

# --- Test Cases sheet: append three new test case rows (TC_005, TC_006, TC_007) ---

# Row 14 - TC_005: blank Email ID, valid Password
$wsCases.Range("A14").Value = "TC_005"
$wsCases.Range("B14").Value = "Login "
$wsCases.Range("C14").Value = "awsomeqa Login Page"
$wsCases.Range("D14").Value = "1. Launch Browser`n2. Go to URL https://awesomeqa.com/ui/`n3. Click on My Account`n4. Click on Login"
$wsCases.Range("E14").Value = "Email ID:- ____________`nPassword:- xyz@1234"
$wsCases.Range("F14").Value = "Verify login with blank Email ID and Valid Password`n"
$wsCases.Range("G14").Value = "P0"
$wsCases.Range("H14").Value = "Invalid Credentials"
$wsCases.Rows.Item(14).RowHeight = 96
$wsCases.Range("D14:H14").WrapText = $true

# Row 15 - TC_006: valid Email ID, blank Password
$wsCases.Range("A15").Value = "TC_006"
$wsCases.Range("B15").Value = "Login "
$wsCases.Range("C15").Value = "awsomeqa Login Page"
$wsCases.Range("D15").Value = "1. Launch Browser`n2. Go to URL https://awesomeqa.com/ui/`n3. Click on My Account`n4. Click on Login"
$wsCases.Range("E15").Value = "Email ID:- xyz@gmail.com`nPassword:- ____________"
$wsCases.Range("F15").Value = "Verify login with valid Email ID`nand blank Password"
$wsCases.Range("G15").Value = "P0"
$wsCases.Range("H15").Value = "Invalid Credentials"
$wsCases.Rows.Item(15).RowHeight = 96.75
$wsCases.Range("D15:H15").WrapText = $true

# Row 16 - TC_007: blank Email ID, blank Password
$wsCases.Range("A16").Value = "TC_007"
$wsCases.Range("B16").Value = "Login "
$wsCases.Range("C16").Value = "awsomeqa Login Page"
$wsCases.Range("D16").Value = "1. Launch Browser`n2. Go to URL https://awesomeqa.com/ui/`n3. Click on My Account`n4. Click on Login"
$wsCases.Range("E16").Value = "Email ID:- ____________`nPassword:- ____________"
$wsCases.Range("F16").Value = "Verify login with blank Email ID`nand Password"
$wsCases.Range("G16").Value = "P0"
$wsCases.Range("H16").Value = "Invalid Credentials"
$wsCases.Rows.Item(16).RowHeight = 94.5
$wsCases.Range("D16:H16").WrapText = $true

# --- Update view/selection state on each sheet to match the saved workbook ---

# "Test Scenarios" tab: selection moves from F4 to F7, no longer the active tab
$wsScenarios.Activate() | Out-Null
$wsScenarios.Range("F7").Select() | Out-Null

# "Sheet3" tab: gains an explicit selection at I18
$wsSheet3.Activate() | Out-Null
$wsSheet3.Range("I18").Select() | Out-Null

# "Test Cases" tab ends up the active / selected tab, with selection at H16
$wsCases.Activate() | Out-Null
$wsCases.Range("H16").Select() | Out-Null
